# Add alejandra.panozo@renewsolutions.com.br to the e-mail distribution
# lists on the "Mails" sheet ("ale panozo agregada a envio de mails").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mails")

# Rows 2-8 (IDs 0-6): append the new address to the "To" column (D),
# which currently only contains maximiliano.esbiza@renewsolutions.com.br
$newTo = "maximiliano.esbiza@renewsolutions.com.br;alejandra.panozo@renewsolutions.com.br"
$toRange = $ws.Range("D2:D8")
$toRange.Value = $newTo
$toRange.WrapText = $false

# Row 9 (ID 7): this mail template is repurposed from "Berry Status Report
# Colombia" to "COL Price List Report y Facturación", its recipients are
# trimmed down and alejandra.panozo@renewsolutions.com.br is added to CC.
$ws.Range("B9").Value = "COL Price List Report y Facturación"
$ws.Range("C9").Value = "Estimados,<br><br>`nSe adjuntan 'Price List Report' y 'XXARX PLANO FACTURACION EXCEL'.<br><br>`nSaludos"
$ws.Range("D9").Value = "EricRodrigues@berryglobal.com;KevinKruger@berryglobal.com;GustavoConforto@berryglobal.com"
$ws.Range("E9").Value = "alejandra.panozo@renewsolutions.com.br;maximiliano.esbiza@renewsolutions.com.br;luisaaranda@berryglobal.com;joaquin.bracci@renewsolutions.com.br"
